# Hortaliza / Papa - weekly price update
# Insert two new weekly observation rows above the existing row 398, pushing
# all subsequent rows down by two (matching the commit's "Fruta / hortaliza,
# semanal" weekly-roll behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 398:399 - everything currently at row 398 onward
# shifts down to row 400 onward.
$ws.Rows("398:399").Insert()

# New row 398
$ws.Cells.Item(398, 1).Value = 8
$ws.Cells.Item(398, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(398, 3).Value = "Coquimbo"
$ws.Cells.Item(398, 4).Value = 44776
$ws.Cells.Item(398, 5).Value = 4
$ws.Cells.Item(398, 6).Value = 100114001
$ws.Cells.Item(398, 7).Value = "Papa"
$ws.Cells.Item(398, 8).Value = "Asterix"
$ws.Cells.Item(398, 9).Value = "1a (cosecha)"
$ws.Cells.Item(398, 10).Value = 2000
$ws.Cells.Item(398, 11).Value = 11500
$ws.Cells.Item(398, 12).Value = 12000
$ws.Cells.Item(398, 13).Value = 11750
$ws.Cells.Item(398, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(398, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(398, 16).Value = 470
$ws.Cells.Item(398, 17).Value = 25
$ws.Cells.Item(398, 18).Value = "Hortaliza"

# New row 399
$ws.Cells.Item(399, 1).Value = 8
$ws.Cells.Item(399, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(399, 3).Value = "Coquimbo"
$ws.Cells.Item(399, 4).Value = 44776
$ws.Cells.Item(399, 5).Value = 4
$ws.Cells.Item(399, 6).Value = 100114001
$ws.Cells.Item(399, 7).Value = "Papa"
$ws.Cells.Item(399, 8).Value = "Cardinal"
$ws.Cells.Item(399, 9).Value = "1a (cosecha)"
$ws.Cells.Item(399, 10).Value = 1600
$ws.Cells.Item(399, 11).Value = 11000
$ws.Cells.Item(399, 12).Value = 11500
$ws.Cells.Item(399, 13).Value = 11250
$ws.Cells.Item(399, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(399, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(399, 16).Value = 450
$ws.Cells.Item(399, 17).Value = 25
$ws.Cells.Item(399, 18).Value = "Hortaliza"
